# scoring_lookup.xlsx: add HAMD6/HRSD6/HRSD17 lookup sheets, rename the
# "forwNames"/"revNames" header columns to "forwItems"/"revItems", and
# update the active-sheet/selection state to match.

$wb = $excel.ActiveWorkbook

# --- Rename the shared header columns on the three existing sheets ---
# (B1 was "forwNames", C1 was "revNames" -> now "forwItems"/"revItems")
foreach ($idx in 1..3) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("B1").Value = "forwItems"
    $ws.Range("C1").Value = "revItems"
}

# --- Add the three new lookup sheets at the end of the workbook ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$hamd6 = $wb.Worksheets.Add([System.Type]::Missing, $last)
$hamd6.Name = "HAMD6total"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$hrsd6 = $wb.Worksheets.Add([System.Type]::Missing, $last)
$hrsd6.Name = "HRSD6total"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$hrsd17 = $wb.Worksheets.Add([System.Type]::Missing, $last)
$hrsd17.Name = "HRSD17total"

# --- HAMD6total (A1:D7) ---
$hamd6.Range("A1").Value = "measName"
$hamd6.Range("B1").Value = "forwItems"
$hamd6.Range("C1").Value = "revItems"
$hamd6.Range("D1").Value = "revInt"
$hamd6.Range("A2").Value = "HAMD6"
$hamd6.Range("B2").Value = "X1"
$hamd6.Range("B3").Value = "X2"
$hamd6.Range("B4").Value = "X7"
$hamd6.Range("B5").Value = "X8"
$hamd6.Range("B6").Value = "X10"
$hamd6.Range("B7").Value = "X13"

# --- HRSD6total (A1:D7) ---
$hrsd6.Range("A1").Value = "measName"
$hrsd6.Range("B1").Value = "forwItems"
$hrsd6.Range("C1").Value = "revItems"
$hrsd6.Range("D1").Value = "revInt"
$hrsd6.Range("A2").Value = "HRSD"
$hrsd6.Range("B2").Value = "X1"
$hrsd6.Range("B3").Value = "X2"
$hrsd6.Range("B4").Value = "X7"
$hrsd6.Range("B5").Value = "X8"
$hrsd6.Range("B6").Value = "X10"
$hrsd6.Range("B7").Value = "X13"

# --- HRSD17total (A1:D18) ---
$hrsd17.Range("A1").Value = "measName"
$hrsd17.Range("B1").Value = "forwItems"
$hrsd17.Range("C1").Value = "revItems"
$hrsd17.Range("D1").Value = "revInt"
$hrsd17.Range("A2").Value = "HRSD"
$hrsd17.Range("B2").Value = "X1"
$hrsd17.Range("B3").Value = "X2"
$hrsd17.Range("B4").Value = "X3"
$hrsd17.Range("B5").Value = "X4"
$hrsd17.Range("B6").Value = "X5"
$hrsd17.Range("B7").Value = "X6"
$hrsd17.Range("B8").Value = "X7"
$hrsd17.Range("B9").Value = "X8"
$hrsd17.Range("B10").Value = "X9"
$hrsd17.Range("B11").Value = "X10"
$hrsd17.Range("B12").Value = "X11"
$hrsd17.Range("B13").Value = "X12"
$hrsd17.Range("B14").Value = "X13"
$hrsd17.Range("B15").Value = "X14"
$hrsd17.Range("B16").Value = "X15"
$hrsd17.Range("B17").Value = "X16"
$hrsd17.Range("B18").Value = "X17"

# --- Per-sheet selections (as left by the editor) ---
$hamd6.Range("B2:B7").Select()

$hrsd6.Range("A1:D2").Select()

$hrsd17.Range("E18").Select()

$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("A1:D1").Select()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("C2:D3").Select()

# Final active sheet/selection: BDItotal, C15
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("C15").Select()
